$d = $word.ActiveDocument

# Locate the title run that currently reads "Registrar Profesional" (inside
# the quoted use-case title near the top of the document) and prefix it with
# the use-case code "CU06 ", keeping it in its own run the way Word does when
# text is typed at that spot, and leaving the auto-managed "_GoBack" bookmark
# wrapped around the last-edited text ("Registrar Profesional").

$range = $d.Content
$found = $range.Find.Execute("Registrar Profesional", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertStart = $range.Start
    $insertEnd = $range.End
    $targetLen = $insertEnd - $insertStart

    # Type the new text right before "Registrar Profesional".
    $insertPoint = $d.Range($insertStart, $insertStart)
    $insertPoint.InsertBefore("CU06 ")

    # Keep "CU06 " from silently re-merging with the preceding run (the
    # opening curly quote) by bracketing it with a throwaway bookmark while
    # we set up the real one.
    $newTextRange = $d.Range($insertStart, $insertStart + 5)
    $d.Bookmarks.Add("_zzTempSplit", $newTextRange)

    # Word automatically drops its "_GoBack" bookmark at the last editing
    # location - re-anchor it around "Registrar Profesional" (this also
    # removes the bookmark from its previous location elsewhere in the doc,
    # since a document only ever has one "_GoBack" bookmark).
    $registrarRange = $d.Range($insertStart + 5, $insertStart + 5 + $targetLen)
    $d.Bookmarks.Add("_GoBack", $registrarRange)

    # Drop the scaffolding bookmark; the run split it enforced persists.
    $d.Bookmarks("_zzTempSplit").Delete()
}
